$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Rename header cells on the existing sheets.
# ---------------------------------------------------------------------
$weekly = $wb.Worksheets.Item("Weekly Quantity")
$weekly.Range("B1").Value = "Weekly_PO_Qty"

$monthly = $wb.Worksheets.Item("Monthly Trend")
$monthly.Range("B1").Value = "Monthly_PO_Qty"

# ---------------------------------------------------------------------
# 2. Add the new "PO Forecast" sheet after the last existing sheet.
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$forecast = $wb.Worksheets.Add($null, $lastSheet)
$forecast.Name = "PO Forecast"

# Header row
$forecast.Range("A1").Value = "ds"
$forecast.Range("B1").Value = "PO_Forecast"
$forecast.Range("C1").Value = "yhat_lower"
$forecast.Range("D1").Value = "yhat_upper"

# Reuse the existing header / date formatting from the "Weekly Quantity"
# sheet so the new sheet matches the workbook's established styling.
[void]$weekly.Range("A1:B1").Copy()
[void]$forecast.Range("A1:D1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

[void]$weekly.Range("A2").Copy()
[void]$forecast.Range("A2:A18").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Data rows
$data = @(
    @(45361.99999999999, 164, 91.75144322501968, 239.6544331433851),
    @(45515.99999999999, 79, 1.524696941170401, 150.9247566639964),
    @(45522.99999999999, 75, 4.471494925372093, 149.4743864481538),
    @(45550.99999999999, 60, -14.95170257079919, 129.3075715836407),
    @(45557.99999999999, 56, -15.79857507181275, 134.3354323222415),
    @(45564.99999999999, 52, -21.10892035912324, 124.8131510328605),
    @(45585.99999999999, 40, -28.92715749907393, 111.9191961980789),
    @(45592.99999999999, 37, -41.65038961145344, 112.3333817235889),
    @(45599.99999999999, 33, -40.20747820955162, 104.652698906299),
    @(45606.99999999999, 29, -46.29481760448503, 102.9127640483353),
    @(45613.99999999999, 25, -53.9566898968576, 97.10047534385792),
    @(45620.99999999999, 21, -51.79803248173664, 91.97434719461111),
    @(45627.99999999999, 17, -57.63753770163204, 89.83909096000171),
    @(45634.99999999999, 13, -54.86249035832387, 91.62845184935874),
    @(45641.99999999999, 9, -66.04362525048423, 82.32742028115224),
    @(45648.99999999999, 6, -74.82676538993645, 72.41174001224417),
    @(45655.99999999999, 2, -69.42586417779945, 77.04647352913857)
)

$row = 2
foreach ($entry in $data) {
    $forecast.Cells.Item($row, 1).Value = $entry[0]
    $forecast.Cells.Item($row, 2).Value = $entry[1]
    $forecast.Cells.Item($row, 3).Value = $entry[2]
    $forecast.Cells.Item($row, 4).Value = $entry[3]
    $row = $row + 1
}

[void]$forecast.Range("A1").Select()
